$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49. This pushes the existing row 49
# (phone "09876543" / birthday 2025-08-12 / points 0) down to row 50,
# and leaves the newly-inserted row 49 blank.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 to match the pattern of the other
# "9876543" customer rows already in the sheet (numeric phone, the
# birthday that used to live on the old row 49, 0 points).
$ws.Cells.Item(49, 1).Value2 = 9876543

# Force the birthday to be stored as literal text (not auto-converted
# to a date serial number) the same way it was on the row above/below.
$ws.Cells.Item(49, 2).NumberFormat = "@"
$ws.Cells.Item(49, 2).Value2 = "2025-08-12"
$ws.Cells.Item(49, 2).ClearFormats()

$ws.Cells.Item(49, 3).Value2 = 0

# The row that got pushed down to 50 keeps its phone text "09876543"
# and its 0 points untouched, but its birthday is cleared out.
$ws.Cells.Item(50, 2).ClearContents()
